$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC2 (rows 24-27), step 4: now becomes the "confirms deletion" scenario
$ws.Range("B27").Value = "Lider de Pessoas confirma a exclusao do Nivel das Competencias"
$ws.Range("D27").Value = "SYSTEM exibe a listagem dos Niveis das Competencias sem o Nivel das Competencias excluido"

# TC3 (rows 34-37), step 4: now becomes the "does NOT confirm deletion" scenario
$ws.Range("B37").Value = "Lider de Pessoas nao confirma a exclusao do Nivel das Competencias"
$ws.Range("D37").Value = "SYSTEM exibe a listagem dos Niveis das Competencias com o Nivel das Competencias nao excluido"

# TC4 (row 45), step 2: typo fix "Niveis das Competencias" -> "Nivel das Competencias"
$ws.Range("B45").Value = "Lider de Pessoas clica na opcao 'Novo' para criar um novo Nivel das Competencias"

# TC6 area (row 66), step 4: keep text as "confirma" (its shared string index moves, text unchanged)
$ws.Range("B66").Value = "Lider de Pessoas confirma a exclusao do Nivel das Competencias"
